$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "None" fingerprint placeholder for rows that don't have one yet,
# matching the formatting (Consolas font, vertical-centered) already used by
# the other fingerprint cells in column C.
$ws.Range("C2").Copy()

$targets = @("C5", "C6", "C9", "C10")
foreach ($addr in $targets) {
    $cell = $ws.Range($addr)
    $cell.Value = "None"
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Update the active selection to reflect where the user ended up.
$ws.Range("C12").Select()
